# Append 6 new data rows (rows 4-9) to "Аркуш1", alternating between two
# sets of values (the second set is the first with " 1" appended), mirroring
# the style (s="9") and layout already used by row 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$setA = @(
    "Тимків Віталій Дмитрович",
    "А0000",
    "Тимків Дмитро Віталійович",
    "Уганда, гасити вагнерів",
    "01.01.2025 Краківець",
    "01.01.2026 Подобовець"
)

$setB = @(
    "Тимків Віталій Дмитрович 1",
    "А0000 1",
    "Тимків Дмитро Віталійович 1",
    "Уганда, гасити вагнерів 1",
    "01.01.2025 Краківець 1",
    "01.01.2026 Подобовець 1"
)

for ($i = 0; $i -lt 6; $i++) {
    $row = 4 + $i

    # Carry over the formatting (style) used on row 3 for columns A-G.
    $ws.Range("A3:G3").Copy($ws.Range("A" + $row + ":G" + $row))

    if (($i % 2) -eq 0) {
        $vals = $setA
    } else {
        $vals = $setB
    }

    $ws.Cells.Item($row, 1).Value = $i + 1
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
